$wb = $excel.ActiveWorkbook

# "Overview" sheet: row for 827e4b46-4ae4-45bb-84a7-74f45a5f979b.md is now ready for handoff
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# "zh-cn" sheet: status + handoff datetime updated
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-03-09 13:27:14"
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-03-09 13:27:14"

# "de-de" sheet: status + handoff datetime updated
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-03-09 13:27:24"
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-03-09 13:27:24"
